$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.340.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.45%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.214.84'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.41%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.39'
$ws.Range('D5').Style = 'Normal'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.28%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.216.55'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.53%  '

# Row 9
$ws.Range('E9').Value = '  -1.81%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.44%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.69'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.07%  '

# Row 12
$ws.Range('E12').Value = '  -3.64%  '

# Row 13
$ws.Range('E13').Value = '  -1.07%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.34%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.738.52'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.35%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.449.85'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.20%  '

# Row 17
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.27'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.34%  '

# Row 18
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.211.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.31%  '

# Row 19
$ws.Range('E19').Value = '  +1.26%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '507.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.29%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.21%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.731'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.42%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.02%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.55%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.14'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.96%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.10%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.49%  '

# Row 28
$ws.Range('E28').Value = '  -2.84%  '

# Row 29
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.133'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +46.62%  '

# Row 30
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.21%  '

# Row 31
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.94'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.65%  '

# Row 32
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.90'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.42%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.56%  '

# Row 35
$ws.Range('E35').Value = '  -5.16%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.43'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.12%  '

# Row 37
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '501.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.68%  '

# Row 38
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.39'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.79%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0772'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.76%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.130'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.03%  '

# Row 41
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0419'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.87%  '

# Row 42
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.85%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.28%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.296'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.89%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.920.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.64%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.71%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.20'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.37%  '

# Row 48
$ws.Range('E48').Value = '  +1.54%  '

# Row 49
$ws.Range('E49').Value = '  -0.04%  '

# Row 50
$ws.Range('E50').Value = '  -1.25%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.45%  '
